# "changes in concise marksheet" - update correct/total marks on the
# marksheet: the "Marking" row's correct-answer count (B11) and the
# "Total" row's correct-answer count (B12) and displayed fraction (E12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking row: correct count 3 -> 5
$ws.Range("B11").Value = 5

# Total row: correct count 54 -> 90
$ws.Range("B12").Value = 90

# Total row: displayed "correct/total" marks fraction
$ws.Range("E12").Value = "90/140"
